# SBM5_PointerTableDocumentation.xlsx edit script
# - Rename sheet1 to "Graphics Pointer Table"
# - Add a new sheet "Menu Text Pointer Table" after it, with a header row copied
#   from sheet1 and a first data row "11A000" in column A
# - Add a note/comment to sheet1!A1 documenting the Pointer Offset column
# - Update selections so the new sheet is the active tab

$wb = $excel.ActiveWorkbook

# --- Rename existing sheet & create the new one -----------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Graphics Pointer Table"

$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Menu Text Pointer Table"

# Copy the header row (values + styles) from the Graphics sheet into the new sheet
$ws1.Range("A1:G1").Copy($ws2.Range("A1"))

# New first data row: pointer offset "11A000"
$ws2.Range("A2").Value = "11A000"

# --- Add reviewer note to the Pointer Offset header on the Graphics sheet ---
$ws1.Range("A1").AddComment("DackR:" + [char]10 + "This is the location of the Pointer. PC Offset.")

# --- Selections / active sheet ----------------------------------------------
$ws1.Range("G2").Select()
$ws2.Activate()
$ws2.Range("B2").Select()
